$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.59%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'31.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.63%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.037"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.55%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07818"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-3.35%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.046"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-21.07%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.785"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.31%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.785"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.91%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9200"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.31%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1748"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-0.43%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07865"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'5.68%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08831"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.67%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03114"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'2.45%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09996"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.13%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001510"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.81%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005852"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-2.73%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'-2.48%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.266"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.86%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3292"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.06%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1282"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-3.95%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.178"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'4.27%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.1810"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'9.85%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04614"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.43%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001243"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.03%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004464"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.38%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001253"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'4.56%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D39").Value = "'0.01741"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-1.73%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04740"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'5.34%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007121"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'4.08%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1355"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'0.34%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002084"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-5.54%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01084"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'10.24%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006050"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-6.28%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'0.38%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.003556"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-59.28%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.8222"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'0.20%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002104"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.38%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002004"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.38%"
$ws.Range("E50").Style = "Normal"
